$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ35009061"
$wb.Worksheets.Item(2).Name = "summ35195914"
$wb.Worksheets.Item(3).Name = "summ35392312"
$wb.Worksheets.Item(4).Name = "summ35644743"
$wb.Worksheets.Item(5).Name = "summ35871949"
$wb.Worksheets.Item(6).Name = "summ36082216"
$wb.Worksheets.Item(7).Name = "summ36307830"
$wb.Worksheets.Item(8).Name = "summ36532786"
$wb.Worksheets.Item(9).Name = "summ36762504"
$wb.Worksheets.Item(10).Name = "summ36983160"
$wb.Worksheets.Item(11).Name = "summ37198843"
$wb.Worksheets.Item(12).Name = "summ37414650"
$wb.Worksheets.Item(13).Name = "summ37613032"
$wb.Worksheets.Item(14).Name = "summ37821687"
$wb.Worksheets.Item(15).Name = "summ38037885"
$wb.Worksheets.Item(16).Name = "summ38341036"
$wb.Worksheets.Item(17).Name = "summ38582055"
$wb.Worksheets.Item(18).Name = "summ38827097"
$wb.Worksheets.Item(19).Name = "summ39096251"
$wb.Worksheets.Item(20).Name = "summ39324993"
$wb.Worksheets.Item(21).Name = "summ39556554"
$wb.Worksheets.Item(22).Name = "summ39777179"
$wb.Worksheets.Item(23).Name = "summ39985855"
$wb.Worksheets.Item(24).Name = "summ40189715"
$wb.Worksheets.Item(25).Name = "summ40398126"
$wb.Worksheets.Item(26).Name = "summ40617525"
$wb.Worksheets.Item(27).Name = "summ40833844"
$wb.Worksheets.Item(28).Name = "summ41064745"
$wb.Worksheets.Item(29).Name = "summ41389081"
$wb.Worksheets.Item(30).Name = "summ41611743"
$wb.Worksheets.Item(31).Name = "summ41839305"
$wb.Worksheets.Item(32).Name = "summ42073004"
$wb.Worksheets.Item(33).Name = "summ42302313"
$wb.Worksheets.Item(34).Name = "summ42523001"
$wb.Worksheets.Item(35).Name = "summ42764417"
$wb.Worksheets.Item(36).Name = "summ43012604"
$wb.Worksheets.Item(37).Name = "summ43235152"
$wb.Worksheets.Item(38).Name = "summ43459438"
$wb.Worksheets.Item(39).Name = "summ43687590"
$wb.Worksheets.Item(40).Name = "summ43924574"
$wb.Worksheets.Item(41).Name = "summ44174848"
$wb.Worksheets.Item(42).Name = "summ44414111"
$wb.Worksheets.Item(43).Name = "summ44670550"
$wb.Worksheets.Item(44).Name = "summ44917565"
$wb.Worksheets.Item(45).Name = "summ45138677"
$wb.Worksheets.Item(46).Name = "summ45360846"
$wb.Worksheets.Item(47).Name = "summ45613837"
$wb.Worksheets.Item(48).Name = "summ45863361"
$wb.Worksheets.Item(49).Name = "summ46125835"
$wb.Worksheets.Item(50).Name = "summ46354628"
